$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 127
$ws.Range("H127").Value = 1152.8485
$ws.Range("I127").Value = 519.1667
$ws.Range("J127").Value = 1913.2667
$ws.Range("K127").Value = 1557.5001
$ws.Range("L127").Value = 5739.800099999999
$ws.Range("M127").Value = 3402.4999
$ws.Range("N127").Value = -15659.8001
# Row 132
$ws.Range("H132").Value = 23814902
$ws.Range("I132").Value = 32263814
$ws.Range("K132").Value = 96791442
$ws.Range("M132").Value = -96788912
# Row 137
$ws.Range("H137").Value = 3743.625
$ws.Range("I137").Value = 3722.652
$ws.Range("K137").Value = 11167.956
$ws.Range("M137").Value = -8617.956
# Row 138
$ws.Range("H138").Value = 3016.58
$ws.Range("I138").Value = 467.76666
$ws.Range("J138").Value = 4108.9287
$ws.Range("K138").Value = 1403.29998
$ws.Range("L138").Value = 12326.7861
$ws.Range("M138").Value = 3736.70002
$ws.Range("N138").Value = -22606.7861

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4076.111
$ws.Range("I32").Value = 4190.5615
$ws.Range("K32").Value = 4190.5615
$ws.Range("M32").Value = -3903.5615
# Row 61
$ws.Range("H61").Value = 1674.6451
$ws.Range("I61").Value = 1255.3334
$ws.Range("K61").Value = 1255.3334
$ws.Range("M61").Value = -1043.3334
# Row 74
$ws.Range("H74").Value = 2895.2
$ws.Range("I74").Value = 3105
$ws.Range("J74").Value = 2265.8
$ws.Range("K74").Value = 3105
$ws.Range("L74").Value = 2265.8
$ws.Range("M74").Value = -2231
$ws.Range("N74").Value = -4013.8
# Row 77
$ws.Range("H77").Value = 2895.2
$ws.Range("I77").Value = 3105
$ws.Range("J77").Value = 2265.8
$ws.Range("K77").Value = 15525
$ws.Range("L77").Value = 11329
$ws.Range("M77").Value = -11157
$ws.Range("N77").Value = -20065
# Row 97
$ws.Range("H97").Value = 716.6
$ws.Range("I97").Value = 775
$ws.Range("J97").Value = 410
$ws.Range("K97").Value = 775
$ws.Range("L97").Value = 410
$ws.Range("M97").Value = -279
$ws.Range("N97").Value = -1402
# Row 132
$ws.Range("H132").Value = 3428.348
$ws.Range("I132").Value = 2204.2307
$ws.Range("K132").Value = 6612.6921
$ws.Range("M132").Value = -4082.6921
# Row 136
$ws.Range("H136").Value = 1674.6451
$ws.Range("I136").Value = 1255.3334
$ws.Range("K136").Value = 3766.0002
$ws.Range("M136").Value = -1216.0002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 42000
$ws.Range("J53").Value = 42000
$ws.Range("L53").Value = 42000
$ws.Range("N53").Value = -43148
# Row 94
$ws.Range("H94").Value = 512.2895
$ws.Range("I94").Value = 569.2593000000001
$ws.Range("J94").Value = 372.45456
$ws.Range("K94").Value = 569.2593000000001
$ws.Range("L94").Value = 372.45456
$ws.Range("M94").Value = -118.2593000000001
$ws.Range("N94").Value = -1274.45456
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 134
$ws.Range("H134").Value = 2758.9048
$ws.Range("I134").Value = 1607
$ws.Range("J134").Value = 4630.75
$ws.Range("K134").Value = 4821
$ws.Range("L134").Value = 13892.25
$ws.Range("M134").Value = -2286
$ws.Range("N134").Value = -18962.25

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6851593.5
$ws.Range("I31").Value = 1356.3541
$ws.Range("K31").Value = 1356.3541
$ws.Range("M31").Value = -1061.3541
# Row 34
$ws.Range("H34").Value = 6851593.5
$ws.Range("I34").Value = 1356.3541
$ws.Range("K34").Value = 1356.3541
$ws.Range("M34").Value = -1154.3541
# Row 99
$ws.Range("H99").Value = 9529399
$ws.Range("I99").Value = 18185982
$ws.Range("J99").Value = 7157
$ws.Range("K99").Value = 18185982
$ws.Range("L99").Value = 7157
$ws.Range("M99").Value = -18184484
$ws.Range("N99").Value = -10153
# Row 126
$ws.Range("H126").Value = 9529399
$ws.Range("I126").Value = 18185982
$ws.Range("J126").Value = 7157
$ws.Range("K126").Value = 54557946
$ws.Range("L126").Value = 21471
$ws.Range("M126").Value = -54555476
$ws.Range("N126").Value = -26411
# Row 132
$ws.Range("H132").Value = 1928.7759
$ws.Range("I132").Value = 1549.06
$ws.Range("K132").Value = 4647.18
$ws.Range("M132").Value = -2117.18
# Row 134
$ws.Range("H134").Value = 3321.2983
$ws.Range("I134").Value = 4287
$ws.Range("K134").Value = 12861
$ws.Range("M134").Value = -10326

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1345.75
$ws.Range("J46").Value = 1626.6666
$ws.Range("L46").Value = 4879.9998
$ws.Range("N46").Value = -5061.9998
# Row 51
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 9000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -8540
$ws.Range("N51").ClearContents()
# Row 57
$ws.Range("H57").Value = 2894.111
$ws.Range("I57").Value = 500
$ws.Range("J57").Value = 3034.9412
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 9104.8236
$ws.Range("M57").Value = -941
$ws.Range("N57").Value = -10222.8236
# Row 58
$ws.Range("H58").Value = 3680
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 4350
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 13050
$ws.Range("M58").Value = -2872
$ws.Range("N58").Value = -13306
# Row 63
$ws.Range("H63").Value = 3317.7144
$ws.Range("I63").Value = 3317.7144
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 9953.143199999999
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -9204.143199999999
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 3317.7144
$ws.Range("I66").Value = 3317.7144
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 29859.4296
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -26115.4296
$ws.Range("N66").ClearContents()
# Row 69
$ws.Range("H69").Value = 2631.7778
$ws.Range("I69").Value = 873
$ws.Range("J69").Value = 4038.8
$ws.Range("K69").Value = 2619
$ws.Range("L69").Value = 12116.4
$ws.Range("M69").Value = -1808
$ws.Range("N69").Value = -13738.4
# Row 72
$ws.Range("H72").Value = 2631.7778
$ws.Range("I72").Value = 873
$ws.Range("J72").Value = 4038.8
$ws.Range("K72").Value = 7857
$ws.Range("L72").Value = 36349.2
$ws.Range("M72").Value = -3801
$ws.Range("N72").Value = -44461.2
# Row 80
$ws.Range("H80").Value = 9466.166999999999
$ws.Range("I80").Value = 7798
$ws.Range("J80").Value = 9799.799999999999
$ws.Range("K80").Value = 23394
$ws.Range("L80").Value = 29399.4
$ws.Range("M80").Value = -22458
$ws.Range("N80").Value = -31271.4
# Row 82
$ws.Range("H82").Value = 3971
$ws.Range("I82").Value = 1006.5
$ws.Range("J82").Value = 9900
$ws.Range("K82").Value = 3019.5
$ws.Range("L82").Value = 29700
$ws.Range("M82").Value = -2613.5
$ws.Range("N82").Value = -30512
# Row 83
$ws.Range("H83").Value = 9466.166999999999
$ws.Range("I83").Value = 7798
$ws.Range("J83").Value = 9799.799999999999
$ws.Range("K83").Value = 70182
$ws.Range("L83").Value = 88198.2
$ws.Range("M83").Value = -65502
$ws.Range("N83").Value = -97558.2
# Row 85
$ws.Range("H85").Value = 3971
$ws.Range("I85").Value = 1006.5
$ws.Range("J85").Value = 9900
$ws.Range("K85").Value = 3019.5
$ws.Range("L85").Value = 29700
$ws.Range("M85").Value = -1615.5
$ws.Range("N85").Value = -32508
# Row 97
$ws.Range("H97").Value = 331.33334
$ws.Range("I97").Value = 262.4
$ws.Range("J97").Value = 417.5
$ws.Range("K97").Value = 787.1999999999999
$ws.Range("L97").Value = 1252.5
$ws.Range("M97").Value = -291.1999999999999
$ws.Range("N97").Value = -2244.5
# Row 100
$ws.Range("H100").Value = 2483.75
$ws.Range("J100").Value = 2483.75
$ws.Range("L100").Value = 7451.25
$ws.Range("N100").Value = -9073.25
# Row 103
$ws.Range("H103").Value = 846.3333
$ws.Range("I103").Value = 410
$ws.Range("K103").Value = 1230
$ws.Range("M103").Value = -351
# Row 131
$ws.Range("H131").Value = 835.6070999999999
$ws.Range("I131").Value = 491.66666
$ws.Range("J131").Value = 929.4091
$ws.Range("K131").Value = 1474.99998
$ws.Range("L131").Value = 2788.2273
$ws.Range("M131").Value = 3565.00002
$ws.Range("N131").Value = -12868.2273
# Row 137
$ws.Range("H137").Value = 3776.182
$ws.Range("I137").Value = 960
$ws.Range("J137").Value = 4832.25
$ws.Range("K137").Value = 2880
$ws.Range("L137").Value = 14496.75
$ws.Range("M137").Value = 2220
$ws.Range("N137").Value = -24696.75

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6759.7607
$ws.Range("I70").Value = 5836.467
$ws.Range("J70").Value = 8490.9375
$ws.Range("K70").Value = 5836.467
$ws.Range("L70").Value = 8490.9375
$ws.Range("M70").Value = -5566.467
$ws.Range("N70").Value = -9030.9375
# Row 73
$ws.Range("H73").Value = 6759.7607
$ws.Range("I73").Value = 5836.467
$ws.Range("J73").Value = 8490.9375
$ws.Range("K73").Value = 5836.467
$ws.Range("L73").Value = 8490.9375
$ws.Range("M73").Value = -4900.467
$ws.Range("N73").Value = -10362.9375
# Row 80
$ws.Range("H80").Value = 13160376
$ws.Range("I80").Value = 20835592
$ws.Range("J80").Value = 2861.4285
$ws.Range("K80").Value = 20835592
$ws.Range("L80").Value = 2861.4285
$ws.Range("M80").Value = -20834594
$ws.Range("N80").Value = -4857.4285
# Row 83
$ws.Range("H83").Value = 13160376
$ws.Range("I83").Value = 20835592
$ws.Range("J83").Value = 2861.4285
$ws.Range("K83").Value = 104177960
$ws.Range("L83").Value = 14307.1425
$ws.Range("M83").Value = -104172968
$ws.Range("N83").Value = -24291.1425
# Row 102
$ws.Range("H102").Value = 1420.3704
$ws.Range("I102").Value = 1060.3611
$ws.Range("J102").Value = 2140.389
$ws.Range("K102").Value = 1060.3611
$ws.Range("L102").Value = 2140.389
$ws.Range("M102").Value = 561.6388999999999
$ws.Range("N102").Value = -5384.389
# Row 132
$ws.Range("H132").Value = 2392.8
$ws.Range("I132").Value = 1474.1154
$ws.Range("J132").Value = 3649.9473
$ws.Range("K132").Value = 4422.3462
$ws.Range("L132").Value = 10949.8419
$ws.Range("M132").Value = -1892.3462
$ws.Range("N132").Value = -16009.8419

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 49
$ws.Range("H49").Value = 25950
$ws.Range("I49").Value = 8000
$ws.Range("K49").Value = 8000
$ws.Range("M49").Value = -7770
# Row 100
$ws.Range("H100").Value = 991.1429000000001
$ws.Range("I100").Value = 739.8333
$ws.Range("J100").Value = 2499
$ws.Range("K100").Value = 1479.6666
$ws.Range("L100").Value = 4998
$ws.Range("M100").Value = -938.6666
$ws.Range("N100").Value = -6080
# Row 132
$ws.Range("H132").Value = 4764652
$ws.Range("I132").Value = 2558.9614
$ws.Range("K132").Value = 7676.8842
$ws.Range("M132").Value = -5146.8842

